# Daily attendance processing - reverse the order of comma-separated
# entries in the "Recorded By" column (G) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ','
        if ($parts.Count -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }
            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }
            $cell.Value2 = [string]::Join(', ', $reversed)
        }
    }
}
